$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.931.69'
$ws.Range("E2").Value = '  -2.19%  '
$ws.Range("D3").Value = '2.568.44'
$ws.Range("E3").Value = '  -2.78%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '515.29'
$ws.Range("E5").Value = '  -2.40%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '139.92'
$ws.Range("E6").Value = '  -3.36%  '
$ws.Range("E7").Value = '  -0.07%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.560'
$ws.Range("E8").Value = '  -1.65%  '
$ws.Range("D9").Value = '2.584.36'
$ws.Range("E9").Value = '  -2.78%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.43'
$ws.Range("E10").Value = '  -3.46%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0995'
$ws.Range("E11").Value = '  -4.31%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.328'
$ws.Range("E12").Value = '  -2.51%  '
$ws.Range("E13").Value = '  +0.31%  '
$ws.Range("D14").Value = '3.021.69'
$ws.Range("E14").Value = '  -2.75%  '
$ws.Range("D15").Value = '57.873.98'
$ws.Range("E15").Value = '  -2.18%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '20.11'
$ws.Range("E16").Value = '  -4.49%  '
$ws.Range("D17").Value = '2.579.46'
$ws.Range("E17").Value = '  -1.73%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0000132'
$ws.Range("E18").Value = '  -3.40%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '334.02'
$ws.Range("E19").Value = '  -2.27%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.28'
$ws.Range("E20").Value = '  -3.80%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.08'
$ws.Range("E21").Value = '  -4.58%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.35'
$ws.Range("E22").Value = '  +0.11%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.998'
$ws.Range("E23").Value = '  +0.02%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '65.51'
$ws.Range("E24").Value = '  +0.38%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.167'
$ws.Range("E25").Value = '  -0.65%  '
$ws.Range("B26").Value = 'Binance-PegBSC-USD'
$ws.Range("C26").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.999'
$ws.Range("E26").Value = '  +0.31%  '
$ws.Range("B27").Value = 'WrappedeETH'
$ws.Range("C27").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D27").Value = '2.706.16'
$ws.Range("E27").Value = '  -1.95%  '
$ws.Range("B28").Value = 'Polygon'
$ws.Range("C28").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.398'
$ws.Range("E28").Value = '  -4.72%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.95'
$ws.Range("E29").Value = '  -4.22%  '
$ws.Range("E30").Value = '  -0.05%  '
$ws.Range("D31").Value = '0.0₃0730'
$ws.Range("E31").Value = '  -8.71%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.01'
$ws.Range("E32").Value = '  -6.91%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.56'
$ws.Range("E33").Value = '  -3.24%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '18.60'
$ws.Range("E34").Value = '  -1.83%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '148.81'
$ws.Range("E35").Value = '  -0.94%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.92'
$ws.Range("E36").Value = '  -6.89%  '
$ws.Range("E37").Value = '  -7.35%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.847'
$ws.Range("E38").Value = '  -2.93%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '35.96'
$ws.Range("E39").Value = '  -1.74%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.45'
$ws.Range("E40").Value = '  -3.27%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.820'
$ws.Range("E41").Value = '  -11.11%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.48'
$ws.Range("E42").Value = '  -5.10%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.997'
$ws.Range("E43").Value = '  -0.05%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '271.61'
$ws.Range("E44").Value = '  -0.24%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.594'
$ws.Range("E45").Value = '  -1.60%  '
$ws.Range("E46").Value = '  +0.41%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0941'
$ws.Range("E47").Value = '  -3.61%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0515'
$ws.Range("E48").Value = '  -4.25%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '18.37'
$ws.Range("E49").Value = '  -5.39%  '
$ws.Range("D50").Value = '1.958.76'
$ws.Range("E50").Value = '  -4.39%  '
$ws.Range("B51").Value = 'VeChain'
$ws.Range("C51").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0218'
$ws.Range("E51").Value = '  -5.05%  '
